$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for a689349f... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-05 06:51:14"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime for a689349f... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-05 06:51:08"
$wsZhCn.Range("K3").Value = "2016-09-05 06:51:54"

# "de-de" sheet: Latest HO Xliff Generate Date / Correspond Handback DateTime for a689349f... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-05 06:51:14"
$wsDeDe.Range("K3").Value = "2016-09-05 06:52:02"
